{"js": "// Office.js (Word JavaScript API) edit script.\n// Replaces the 100 \"AxB=\" multiplication prompts in the single 5-column\n// table (20 rows x 5 cols, row-major order matches document order) with\n// their updated values, per the target diff. Some \"old\" values repeat\n// (e.g. \"75\u00d738=\" appears twice with different replacements), so the\n// edit addresses each table cell by (row, col) position rather than doing\n// a global text find/replace.\n\nconst OLD_VALUES = [\"44\u00d791=\", \"71\u00d767=\", \"69\u00d761=\", \"92\u00d729=\", \"14\u00d725=\", \"74\u00d727=\", \"63\u00d793=\", \"76\u00d762=\", \"18\u00d747=\", \"30\u00d749=\", \"13\u00d767=\", \"70\u00d735=\", \"66\u00d766=\", \"50\u00d757=\", \"98\u00d734=\", \"55\u00d714=\", \"75\u00d738=\", \"39\u00d737=\", \"19\u00d782=\", \"50\u00d727=\", \"11\u00d757=\", \"75\u00d738=\", \"64\u00d783=\", \"44\u00d713=\", \"99\u00d754=\", \"43\u00d794=\", \"40\u00d770=\", \"94\u00d787=\", \"18\u00d723=\", \"47\u00d750=\", \"19\u00d744=\", \"66\u00d784=\", \"61\u00d764=\", \"36\u00d771=\", \"51\u00d765=\", \"100\u00d782=\", \"74\u00d743=\", \"86\u00d799=\", \"99\u00d767=\", \"15\u00d771=\", \"95\u00d744=\", \"81\u00d744=\", \"99\u00d739=\", \"48\u00d712=\", \"15\u00d765=\", \"84\u00d763=\", \"49\u00d752=\", \"41\u00d724=\", \"52\u00d797=\", \"95\u00d784=\", \"85\u00d721=\", \"59\u00d727=\", \"46\u00d751=\", \"60\u00d734=\", \"88\u00d786=\", \"57\u00d774=\", \"61\u00d768=\", \"28\u00d780=\", \"10\u00d759=\", \"52\u00d756=\", \"56\u00d723=\", \"68\u00d729=\", \"97\u00d743=\", \"38\u00d743=\", \"56\u00d793=\", \"37\u00d771=\", \"80\u00d772=\", \"59\u00d738=\", \"81\u00d757=\", \"59\u00d776=\", \"15\u00d745=\", \"12\u00d771=\", \"31\u00d776=\", \"85\u00d767=\", \"87\u00d749=\", \"88\u00d779=\", \"27\u00d716=\", \"63\u00d799=\", \"75\u00d790=\", \"16\u00d739=\", \"98\u00d729=\", \"49\u00d787=\", \"41\u00d726=\", \"28\u00d742=\", \"75\u00d792=\", \"12\u00d750=\", \"59\u00d747=\", \"39\u00d756=\", \"41\u00d710=\", \"62\u00d714=\", \"80\u00d782=\", \"77\u00d797=\", \"27\u00d761=\", \"75\u00d798=\", \"70\u00d7100=\", \"81\u00d797=\", \"17\u00d721=\", \"73\u00d746=\", \"50\u00d750=\", \"44\u00d770=\"];\nconst NEW_VALUES = [\"64\u00d752=\", \"42\u00d749=\", \"11\u00d717=\", \"49\u00d722=\", \"16\u00d728=\", \"51\u00d747=\", \"18\u00d736=\", \"74\u00d715=\", \"47\u00d717=\", \"53\u00d736=\", \"74\u00d756=\", \"92\u00d762=\", \"44\u00d732=\", \"31\u00d754=\", \"43\u00d750=\", \"89\u00d761=\", \"75\u00d755=\", \"33\u00d7100=\", \"100\u00d765=\", \"79\u00d793=\", \"93\u00d799=\", \"47\u00d725=\", \"58\u00d719=\", \"41\u00d733=\", \"65\u00d711=\", \"79\u00d744=\", \"66\u00d763=\", \"73\u00d738=\", \"98\u00d766=\", \"19\u00d718=\", \"76\u00d765=\", \"75\u00d777=\", \"29\u00d770=\", \"23\u00d769=\", \"98\u00d794=\", \"13\u00d716=\", \"83\u00d718=\", \"90\u00d759=\", \"42\u00d771=\", \"42\u00d774=\", \"25\u00d735=\", \"18\u00d752=\", \"58\u00d735=\", \"81\u00d782=\", \"76\u00d768=\", \"60\u00d737=\", \"48\u00d795=\", \"43\u00d721=\", \"35\u00d792=\", \"58\u00d742=\", \"64\u00d725=\", \"70\u00d726=\", \"13\u00d7100=\", \"43\u00d769=\", \"79\u00d752=\", \"73\u00d743=\", \"21\u00d738=\", \"40\u00d745=\", \"27\u00d728=\", \"16\u00d772=\", \"18\u00d714=\", \"80\u00d786=\", \"24\u00d722=\", \"15\u00d746=\", \"92\u00d787=\", \"67\u00d771=\", \"56\u00d744=\", \"30\u00d777=\", \"89\u00d790=\", \"25\u00d713=\", \"58\u00d782=\", \"63\u00d737=\", \"23\u00d730=\", \"51\u00d732=\", \"39\u00d759=\", \"46\u00d763=\", \"86\u00d725=\", \"84\u00d775=\", \"61\u00d746=\", \"87\u00d723=\", \"59\u00d767=\", \"22\u00d745=\", \"25\u00d786=\", \"71\u00d771=\", \"13\u00d752=\", \"98\u00d753=\", \"67\u00d750=\", \"63\u00d780=\", \"26\u00d722=\", \"17\u00d754=\", \"14\u00d768=\", \"35\u00d723=\", \"84\u00d776=\", \"40\u00d779=\", \"39\u00d792=\", \"33\u00d737=\", \"25\u00d746=\", \"28\u00d775=\", \"26\u00d717=\", \"56\u00d791=\"];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst COLS = 5;\nconst rows = table.rowCount;\n\n// Collect cell proxies in row-major order and load their current text.\nconst cells = [];\nfor (let r = 0; r < rows; r++) {\n  for (let c = 0; c < COLS; c++) {\n    const cell = table.getCell(r, c);\n    cell.load(\"value\");\n    cells.push(cell);\n  }\n}\nawait context.sync();\n\nconst count = Math.min(cells.length, NEW_VALUES.length);\nfor (let i = 0; i < count; i++) {\n  const cell = cells[i];\n  const expected = OLD_VALUES[i];\n  const current = (cell.value || \"\").trim();\n  // Only overwrite the cell if it still holds the value we expect to\n  // replace (defensive \u2014 keeps this idempotent / safe to re-run).\n  if (expected === undefined || current === expected) {\n    cell.value = NEW_VALUES[i];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Replaces the 100 \"AxB=\" multiplication prompts in the single 5-column\n# table (20 rows x 5 cols, row-major order matches document order) with\n# their updated values, per the target diff. Some \"old\" values repeat\n# (e.g. \"75\u00d738=\" appears twice with different replacements), so the\n# edit addresses each table cell by (row, col) position via Cell(r, c)\n# rather than doing a global Find/Replace.\n\n$OldValues = @(\"44\u00d791=\",\"71\u00d767=\",\"69\u00d761=\",\"92\u00d729=\",\"14\u00d725=\",\"74\u00d727=\",\"63\u00d793=\",\"76\u00d762=\",\"18\u00d747=\",\"30\u00d749=\",\"13\u00d767=\",\"70\u00d735=\",\"66\u00d766=\",\"50\u00d757=\",\"98\u00d734=\",\"55\u00d714=\",\"75\u00d738=\",\"39\u00d737=\",\"19\u00d782=\",\"50\u00d727=\",\"11\u00d757=\",\"75\u00d738=\",\"64\u00d783=\",\"44\u00d713=\",\"99\u00d754=\",\"43\u00d794=\",\"40\u00d770=\",\"94\u00d787=\",\"18\u00d723=\",\"47\u00d750=\",\"19\u00d744=\",\"66\u00d784=\",\"61\u00d764=\",\"36\u00d771=\",\"51\u00d765=\",\"100\u00d782=\",\"74\u00d743=\",\"86\u00d799=\",\"99\u00d767=\",\"15\u00d771=\",\"95\u00d744=\",\"81\u00d744=\",\"99\u00d739=\",\"48\u00d712=\",\"15\u00d765=\",\"84\u00d763=\",\"49\u00d752=\",\"41\u00d724=\",\"52\u00d797=\",\"95\u00d784=\",\"85\u00d721=\",\"59\u00d727=\",\"46\u00d751=\",\"60\u00d734=\",\"88\u00d786=\",\"57\u00d774=\",\"61\u00d768=\",\"28\u00d780=\",\"10\u00d759=\",\"52\u00d756=\",\"56\u00d723=\",\"68\u00d729=\",\"97\u00d743=\",\"38\u00d743=\",\"56\u00d793=\",\"37\u00d771=\",\"80\u00d772=\",\"59\u00d738=\",\"81\u00d757=\",\"59\u00d776=\",\"15\u00d745=\",\"12\u00d771=\",\"31\u00d776=\",\"85\u00d767=\",\"87\u00d749=\",\"88\u00d779=\",\"27\u00d716=\",\"63\u00d799=\",\"75\u00d790=\",\"16\u00d739=\",\"98\u00d729=\",\"49\u00d787=\",\"41\u00d726=\",\"28\u00d742=\",\"75\u00d792=\",\"12\u00d750=\",\"59\u00d747=\",\"39\u00d756=\",\"41\u00d710=\",\"62\u00d714=\",\"80\u00d782=\",\"77\u00d797=\",\"27\u00d761=\",\"75\u00d798=\",\"70\u00d7100=\",\"81\u00d797=\",\"17\u00d721=\",\"73\u00d746=\",\"50\u00d750=\",\"44\u00d770=\")\n$NewValues = @(\"64\u00d752=\",\"42\u00d749=\",\"11\u00d717=\",\"49\u00d722=\",\"16\u00d728=\",\"51\u00d747=\",\"18\u00d736=\",\"74\u00d715=\",\"47\u00d717=\",\"53\u00d736=\",\"74\u00d756=\",\"92\u00d762=\",\"44\u00d732=\",\"31\u00d754=\",\"43\u00d750=\",\"89\u00d761=\",\"75\u00d755=\",\"33\u00d7100=\",\"100\u00d765=\",\"79\u00d793=\",\"93\u00d799=\",\"47\u00d725=\",\"58\u00d719=\",\"41\u00d733=\",\"65\u00d711=\",\"79\u00d744=\",\"66\u00d763=\",\"73\u00d738=\",\"98\u00d766=\",\"19\u00d718=\",\"76\u00d765=\",\"75\u00d777=\",\"29\u00d770=\",\"23\u00d769=\",\"98\u00d794=\",\"13\u00d716=\",\"83\u00d718=\",\"90\u00d759=\",\"42\u00d771=\",\"42\u00d774=\",\"25\u00d735=\",\"18\u00d752=\",\"58\u00d735=\",\"81\u00d782=\",\"76\u00d768=\",\"60\u00d737=\",\"48\u00d795=\",\"43\u00d721=\",\"35\u00d792=\",\"58\u00d742=\",\"64\u00d725=\",\"70\u00d726=\",\"13\u00d7100=\",\"43\u00d769=\",\"79\u00d752=\",\"73\u00d743=\",\"21\u00d738=\",\"40\u00d745=\",\"27\u00d728=\",\"16\u00d772=\",\"18\u00d714=\",\"80\u00d786=\",\"24\u00d722=\",\"15\u00d746=\",\"92\u00d787=\",\"67\u00d771=\",\"56\u00d744=\",\"30\u00d777=\",\"89\u00d790=\",\"25\u00d713=\",\"58\u00d782=\",\"63\u00d737=\",\"23\u00d730=\",\"51\u00d732=\",\"39\u00d759=\",\"46\u00d763=\",\"86\u00d725=\",\"84\u00d775=\",\"61\u00d746=\",\"87\u00d723=\",\"59\u00d767=\",\"22\u00d745=\",\"25\u00d786=\",\"71\u00d771=\",\"13\u00d752=\",\"98\u00d753=\",\"67\u00d750=\",\"63\u00d780=\",\"26\u00d722=\",\"17\u00d754=\",\"14\u00d768=\",\"35\u00d723=\",\"84\u00d776=\",\"40\u00d779=\",\"39\u00d792=\",\"33\u00d737=\",\"25\u00d746=\",\"28\u00d775=\",\"26\u00d717=\",\"56\u00d791=\")\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        if ($i -ge $NewValues.Length) { break }\n        $cell = $t.Cell($r, $c)\n        $rng = $cell.Range\n        $current = $rng.Text\n        $expected = $OldValues[$i] + \"`r`a\"\n        # Only overwrite the cell if it still holds the value we expect to\n        # replace (defensive -- keeps this idempotent / safe to re-run).\n        if ($current -eq $expected) {\n            $rng.Text = $NewValues[$i]\n        }\n        $i++\n    }\n}\n"}
